$wb = $excel.ActiveWorkbook

# The localization status for two files moved from "Ready for handoff" to
# "In Translation": 5bc84813-9d30-4b01-b73b-5571ae8e79d4.md (row 3) and
# 61c4d048-704d-4d88-8738-371cb37d0720.md (row 4). Update the per-language
# Status column on each language sheet, as well as the summary columns on
# the Overview sheet, for those two rows. The third file,
# 6755b1d2-6be2-4a52-b1d1-d0fecea89a73.md (row 5), keeps "Ready for handoff".

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(3, 5).Value = "In Translation"
$wsOverview.Cells.Item(3, 6).Value = "In Translation"
$wsOverview.Cells.Item(4, 5).Value = "In Translation"
$wsOverview.Cells.Item(4, 6).Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Item(3, 3).Value = "In Translation"
$wsZhCn.Cells.Item(4, 3).Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Item(3, 3).Value = "In Translation"
$wsDeDe.Cells.Item(4, 3).Value = "In Translation"
